$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row for 2022-Q3 right after the
#    header row, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.21

# Row 2's "A" cell should carry the same style as the other index cells
# below it (style index copied from A3, which already has it).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# Renumber the sequential index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# 2. Add a new "2022-Q3" sheet, positioned right after "总计" (i.e. before
#    the current "2022-Q2" sheet). Clone the existing "2022-Q2" sheet so the
#    layout/styles match the other per-quarter sheets, then overwrite its
#    contents with the 2022-Q3 fund holdings.
# ---------------------------------------------------------------------------
$q2old = $wb.Worksheets.Item(2)
$q2old.Copy($q2old)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Drop the extra data rows inherited from the copied "2022-Q2" sheet (it has
# 10 data rows, the new sheet only needs 5).
$q3.Range("A7:H11").Delete()

$data = @(
    @(0, "002810", "金信转型创新成长灵活配置混合", "4.06", "93.22", "5.06", "0.2054", 5),
    @(1, "002068", "东方多策略灵活配置混合C", "0.26", "55.14", "2.47", "0.0064", 6),
    @(2, "004402", "金信民旺债券C", "0.09", "23.60", "1.14", "0.0010", 9),
    @(3, "004222", "金信民旺债券A", "0.08", "23.60", "1.14", "0.0009", 9),
    @(4, "400023", "东方多策略灵活配置混合A", "0.03", "55.14", "2.47", "0.0007", 6)
)

$r = 2
foreach ($row in $data) {
    $q3.Range("A$r").Value = $row[0]

    # Fund code (B) and the regime-dependent numeric-looking text columns
    # (D-G) must stay text (e.g. "002810" must keep its leading zero), so
    # force text format before assigning, then strip the style again so the
    # cell ends up with no explicit style index (matching the other rows).
    foreach ($col in @("B", "D", "E", "F", "G")) {
        $idx = switch ($col) { "B" {1} "D" {3} "E" {4} "F" {5} "G" {6} }
        $cell = $q3.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value = $row[$idx]
        $cell.ClearFormats()
    }

    $q3.Range("C$r").Value = $row[2]
    $q3.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# Keep the originally-active sheet selected, matching the pre-edit workbook.
$wb.Worksheets.Item(1).Activate()
